# Trade #81 closed at 2026-02-17 08:59:00 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.5
$summary.Range("B4").Value = 0.51
$summary.Range("B5").Value = 0.13
$summary.Range("B6").Value = 81
$summary.Range("B7").Value = 34
$summary.Range("B9").Value = 41.98

# ---- Strategy Status sheet (MarketMaking row) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.5
$status.Range("D4").Value = 81
$status.Range("E4").Value = 0.51
$status.Range("F4").Value = 0.5
$status.Range("G4").Value = 41.98

# ---- New trade row (#81) appended to "All Trades" and "MarketMaking" sheets ----
$newRow = @(
    81,
    "2026-02-17",
    "08:58:54",
    "MarketMaking",
    "DOWN",
    0.92,
    0.9399999999999999,
    "CLOSED",
    2.1739,
    0.02,
    100.5,
    0,
    0,
    0.6,
    "Normal spread capture: 19600 bps",
    "early_exit",
    0.13
)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 82
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($col -eq 2) {
            # Force the date-like string to stay plain text instead of being
            # auto-converted to a date serial number (matches columns B
            # elsewhere in the sheet, stored as literal text).
            $cell.NumberFormat = "@"
            $cell.Value = $newRow[$i]
            $cell.ClearFormats()
        } else {
            $cell.Value = $newRow[$i]
        }
    }
}
